$p = $ppt.ActivePresentation

# 1. Remove the old "Worked out example" slide (position 7).
$p.Slides.Item(7).Delete()

# 2. Tweak the body text on the "Compiling to the dRMT architecture" slide
#    (now shifted up to position 7) to add "e.g., " before "match capacity".
$compSlide = $p.Slides.Item(7)
$contentShape = $compSlide.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)
$para.Text = "Generate a periodic schedule of operations across processors without violating resource constraints (e.g., match capacity, action capacity, and memory capacity)"

# 3. Add a new slide before "Future work" (currently position 14) for the
#    new "Hardware costs of dRMT" content, using the same "Title and
#    Content" layout used by the surrounding slides.
$newSlide = $p.Slides.Add(14, 2)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Hardware costs of dRMT"
